$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Helper: write a literal text value into a cell without letting Excel's
# autocorrect reinterpret date-looking strings (e.g. "08.10.2022") as a
# real date serial (which would also mint a brand-new number-format style).
# We stage the literal string as a text-formula result, copy it, and paste
# only the value into the destination; Paste-Values does not re-run the
# keyboard-entry date heuristics, so the text and the original General
# style both survive untouched.
function Set-TextValue($cell, [string]$text) {
    $staging = $ws.Range("ZZ1")
    $escaped = $text.Replace("""", """""")
    $staging.Formula = "=""" + $escaped + """"
    $staging.Copy()
    $cell.PasteSpecial(-4163) | Out-Null   # xlPasteValues
    $excel.CutCopyMode = $false
    $staging.Clear()
}

# The summary formula currently lives in row 15 (D15). Push it down to row
# 19 to make room for four new rows (15-18) of timeline data. Inserting
# whole rows copies the row-above formatting down (so column A/B keep
# their existing styles s="2"/s="3" automatically).
$ws.Rows("15:18").Insert()

# --- Row 15: 08.10.2022 -----------------------------------------------
Set-TextValue $ws.Range("A15") "08.10.2022"
$ws.Range("B15").Value = 0.66666666666666663
Set-TextValue $ws.Range("C15") "Project Introduction"
$ws.Range("D15").Value = 60
Set-TextValue $ws.Range("E15") "Outline of Specifications and Terminology: Chord, Riff, Monophone/Polyphone"

# --- Row 17: continuation entry (filled in before row 16, matching the
#     original authoring order reflected by the sharedStrings sequence) --
Set-TextValue $ws.Range("C17") "Wiring Experiment"
Set-TextValue $ws.Range("E17") "Create an Arduino 4X4 Numpad"

# --- Row 16: 09.10.2022 -------------------------------------------------
Set-TextValue $ws.Range("A16") "09.10.2022"
$ws.Range("B16").Value = 0.54166666666666663
Set-TextValue $ws.Range("C16") "Research"
$ws.Range("D16").Value = 190
Set-TextValue $ws.Range("E16") "Existing Technologies, Guitar Hero, MI Digital Guitar, RockSmith and their specifications and Comparision Section"

# Row 18 stays blank (A18/B18 already carry the right style from Insert).

# Update the view to match the saved scroll position/selection.
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("D17").Select()

$wb.Save()
